# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet gains three new trailing columns:
#   H = date              (the report date, 2011-11-15, as text)
#   I = legislator_name   (張慶忠)
#   J = legislator_id     (1347)
# populated for the header row and every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "張慶忠"
$legislatorId = 1347
$reportDate = "2011-11-15"

# Find the last used row/column on the sheet (header row 1 + 23 data rows -> 24).
$lastRow = $ws.Cells.SpecialCells(11).Row
$lastCol = $ws.Cells.SpecialCells(11).Column

# --- Header row -------------------------------------------------------
$ws.Cells.Item(1, $lastCol + 1).Value = "date"
$ws.Cells.Item(1, $lastCol + 2).Value = "legislator_name"
$ws.Cells.Item(1, $lastCol + 3).Value = "legislator_id"

# Copy the header formatting (bold + border + centered) onto the new headers.
$ws.Range($ws.Cells.Item(1, $lastCol), $ws.Cells.Item(1, $lastCol)).Copy()
$ws.Range($ws.Cells.Item(1, $lastCol + 1), $ws.Cells.Item(1, $lastCol + 3)).PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    # Match the plain data-row formatting used by the rest of the table.
    $ws.Range($ws.Cells.Item($r, $lastCol), $ws.Cells.Item($r, $lastCol)).Copy()
    $ws.Range($ws.Cells.Item($r, $lastCol + 1), $ws.Cells.Item($r, $lastCol + 3)).PasteSpecial(-4122)

    # Force the date into the sheet as literal text (not an auto-parsed date serial).
    $ws.Cells.Item($r, $lastCol + 1).NumberFormat = "@"
    $ws.Cells.Item($r, $lastCol + 1).Value = $reportDate

    $ws.Cells.Item($r, $lastCol + 2).Value = $legislatorName
    $ws.Cells.Item($r, $lastCol + 3).Value = $legislatorId
}

$excel.CutCopyMode = $false
